$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 3402.533333333333, 3671, 3175, 0.05369596481323242),
    @(1, 3593.833333333333, 3855, 3257, 0.05112924575805664),
    @(2, 3688.033333333333, 3963, 3381, 0.05392231146494548),
    @(3, 3503.7, 3755, 3143, 0.05264371236165365),
    @(4, 2647.833333333333, 2933, 2307, 0.0554495890935262),
    @(5, 2830.033333333333, 3060, 2460, 0.05255053043365478),
    @(6, 3603.9, 3888, 3214, 0.05511669317881266),
    @(7, 3156.3, 3438, 2789, 0.05420364538828532),
    @(8, 3470.733333333333, 3720, 3117, 0.05325449307759603),
    @(9, 3068.066666666667, 3354, 2789, 0.05127132733662923)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
